$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.795.12"
$ws.Range("E2").Value = "  +5.07%  "
$ws.Range("D3").Value = "2.723.56"
$ws.Range("E3").Value = "  +2.88%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'575.84"
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("D6").Value = "'154.37"
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.40%  "
$ws.Range("D9").Value = "2.747.19"
$ws.Range("E9").Value = "  +3.00%  "
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("D11").Value = "'0.113"
$ws.Range("E11").Value = "  +5.41%  "
$ws.Range("E12").Value = "  +5.53%  "
$ws.Range("E13").Value = "  +3.56%  "
$ws.Range("D14").Value = "3.207.39"
$ws.Range("E14").Value = "  +2.81%  "
$ws.Range("D15").Value = "'26.28"
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("D16").Value = "63.670.95"
$ws.Range("E16").Value = "  +4.89%  "
$ws.Range("E17").Value = "  +6.32%  "
$ws.Range("D18").Value = "2.739.58"
$ws.Range("E18").Value = "  +2.94%  "
$ws.Range("E19").Value = "  +3.35%  "
$ws.Range("E20").Value = "  +2.70%  "
$ws.Range("D21").Value = "'360.19"
$ws.Range("E21").Value = "  +2.73%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'0.996"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("D24").Value = "'0.534"
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("D25").Value = "'65.92"
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("E26").Value = "  +4.22%  "
$ws.Range("E27").Value = "  +4.21%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  +11.76%  "
$ws.Range("D30").Value = "'1.99"
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("D31").Value = "'7.11"
$ws.Range("E31").Value = "  +6.54%  "
$ws.Range("D32").Value = "'171.56"
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("E33").Value = "  +13.41%  "
$ws.Range("D35").Value = "'20.46"
$ws.Range("E35").Value = "  +3.76%  "
$ws.Range("E36").Value = "  +6.95%  "
$ws.Range("E37").Value = "  +9.19%  "
$ws.Range("E38").Value = "  +9.68%  "
$ws.Range("E39").Value = "  +13.36%  "
$ws.Range("D40").Value = "'344.45"
$ws.Range("E40").Value = "  +4.01%  "
$ws.Range("D41").Value = "'4.21"
$ws.Range("E41").Value = "  +5.44%  "
$ws.Range("E42").Value = "  +2.58%  "
$ws.Range("D43").Value = "'5.59"
$ws.Range("E43").Value = "  +7.26%  "
$ws.Range("D44").Value = "'21.68"
$ws.Range("E44").Value = "  +6.64%  "
$ws.Range("D45").Value = "'21.74"
$ws.Range("E45").Value = "  +3.91%  "
$ws.Range("D46").Value = "'0.0590"
$ws.Range("E46").Value = "  +5.37%  "
$ws.Range("D47").Value = "'139.05"
$ws.Range("E47").Value = "  +3.97%  "
$ws.Range("E48").Value = "  +4.47%  "
$ws.Range("D49").Value = "'0.0254"
$ws.Range("E49").Value = "  +3.17%  "
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").Value = "'0.997"
$ws.Range("E51").Value = "  -0.06%  "
